$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row 1: "100" -> "0M"
$t.Cell(1, 1).Range.Text = "0M"

# Row 2: "0.01" -> "0M"
$t.Cell(2, 1).Range.Text = "0M"

# Row 3: "344" -> "0M"
$t.Cell(3, 1).Range.Text = "0M"

# Row 4: "3" -> "203"
$t.Cell(4, 1).Range.Text = "203"

# Row 5: "0.00004" -> "0.00002"
$t.Cell(5, 1).Range.Text = "0.00002"

# Row 6: "0.00006" -> "0.00021"
$t.Cell(6, 1).Range.Text = "0.00021"

# Row 10: "0.00004" -> "0.00005"
$t.Cell(10, 1).Range.Text = "0.00005"

# Row 12: "0.00015" -> "0.00914"
$t.Cell(12, 1).Range.Text = "0.00914"

# Row 44: multi-value tab-separated row collapses to single value "100"
$t.Cell(44, 1).Range.Text = "100"

# Row 45: multi-value tab-separated row collapses to single value "0.01"
$t.Cell(45, 1).Range.Text = "0.01"

# Row 46: multi-value tab-separated row collapses to single value "344"
$t.Cell(46, 1).Range.Text = "344"
